$d = $word.ActiveDocument

# Update the date heading at the top of the document.
$d.Content.Find.Execute("2024-04-13 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-14 Sunday", 2)

# Update the division problems in the table. Each populated row is
# addressed by its (row, column) position rather than by matching old
# text, because a couple of the new values collide with old values
# found elsewhere in the table (e.g. "40÷8=" is both replaced and
# introduced), so a global text find/replace could double-apply.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "50÷3="
$t.Cell(1, 2).Range.Text  = "20÷7="
$t.Cell(1, 3).Range.Text  = "13÷2="
$t.Cell(1, 4).Range.Text  = "68÷3="
$t.Cell(1, 5).Range.Text  = "40÷8="

$t.Cell(5, 1).Range.Text  = "95÷5="
$t.Cell(5, 2).Range.Text  = "90÷8="
$t.Cell(5, 3).Range.Text  = "88÷7="
$t.Cell(5, 4).Range.Text  = "37÷5="
$t.Cell(5, 5).Range.Text  = "22÷3="

$t.Cell(9, 1).Range.Text  = "52÷9="
$t.Cell(9, 2).Range.Text  = "12÷7="
$t.Cell(9, 3).Range.Text  = "77÷9="
$t.Cell(9, 4).Range.Text  = "30÷4="
$t.Cell(9, 5).Range.Text  = "75÷3="

$t.Cell(13, 1).Range.Text = "25÷4="
$t.Cell(13, 2).Range.Text = "22÷2="
$t.Cell(13, 3).Range.Text = "98÷7="
$t.Cell(13, 4).Range.Text = "99÷6="
$t.Cell(13, 5).Range.Text = "42÷3="

$t.Cell(17, 1).Range.Text = "25÷5="
$t.Cell(17, 2).Range.Text = "99÷3="
$t.Cell(17, 3).Range.Text = "96÷5="
$t.Cell(17, 4).Range.Text = "58÷3="
$t.Cell(17, 5).Range.Text = "59÷9="

Write-Host "Done applying updates"
